$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "questMode" column before the existing "questType" column (E)
$ws.Range("E1").EntireColumn.Insert()
$ws.Range("E1").Value = "questMode"

# Add the new "marks" column (H) header
$ws.Range("H1").Value = "marks"

# Row 2: paperCode/instCode become the subjective-question codes; add
# questMode/questType values and the question count + new marks value
$ws.Range("A2").Value = "EEF305"
$ws.Range("B2").Value = "GPK"
$ws.Range("E2").Value = "O"
$ws.Range("F2").Value = "R"
$ws.Range("G2").Value = 38
$ws.Range("H2").Value = 1

# Row 3: same paperCode/instCode, different questMode + question count
$ws.Range("A3").Value = "EEF305"
$ws.Range("B3").Value = "GPK"
$ws.Range("E3").Value = "S"
$ws.Range("F3").Value = "R"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 1

# Header row + data rows get centered alignment
$headerRange = $ws.Range("A1:H1")
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108

$dataRange = $ws.Range("B2:H3")
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108

# Shrink the print scale slightly to fit the new column
$ws.PageSetup.Zoom = 92

# Move the active selection to the new marks column
$ws.Range("H7").Select()
